$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "picture" column (E) file names: old "images/<name>.jpeg" files were
# re-uploaded as "images/<name>.jpg" (and "blueband" became "blue_band").
# Writing in this particular order reproduces the shared-string table order
# produced by the re-upload (new strings are appended in first-use order).
$ws.Range("E3").Value  = "images/margarin_filma.jpg"
$ws.Range("E5").Value  = "images/roti_tawar.jpg"
$ws.Range("E4").Value  = "images/telur_ayam.jpg"
$ws.Range("E2").Value  = "images/blue_band.jpg"
$ws.Range("E11").Value = "images/margarin_filma.jpg"
$ws.Range("E13").Value = "images/roti_tawar.jpg"
$ws.Range("E12").Value = "images/telur_ayam.jpg"
$ws.Range("E10").Value = "images/blue_band.jpg"

# Updated prices for the Alfamart rows (10-13).
$ws.Range("C10").Value = 12500
$ws.Range("C11").Value = 11000
$ws.Range("C12").Value = 25000
$ws.Range("C13").Value = 27000

# Widen column B to fit the content better (Excel quantizes ColumnWidth to
# whole pixels, so 20.6 characters is the closest settable value to the
# ~21.57 stored width).
$ws.Columns("B").ColumnWidth = 20.6

# Move the active selection to E10.
$ws.Range("E10").Activate()
